$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Locator Type" column (D) with header and per-row CSS/Xpath classification
$ws.Range("D1").Value = "Locator Type"
$ws.Range("D2").Value = "CSS"
$ws.Range("D3").Value = "CSS"
$ws.Range("D4").Value = "CSS"
$ws.Range("D5").Value = "CSS"
$ws.Range("D6").Value = "CSS"
$ws.Range("D7").Value = "CSS"
$ws.Range("D8").Value = "CSS"
$ws.Range("D9").Value = "Xpath"
$ws.Range("D10").Value = "Xpath"
$ws.Range("D11").Value = "CSS"
$ws.Range("D12").Value = "CSS"
$ws.Range("D13").Value = "CSS"
$ws.Range("D14").Value = "Xpath"

# Match the persisted column width for column D (~10.5 characters, best-fit)
$ws.Columns.Item(4).ColumnWidth = 9.67

# Move the active selection to D14, matching the saved cursor position
$ws.Range("D14").Select() | Out-Null
